$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 84
$ws.Cells.Item(84, 1).Value = 'Individ 6'
$ws.Cells.Item(84, 2).Value = 15
$ws.Cells.Item(84, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(84, 4).Value = ''
$ws.Cells.Item(84, 5).Value = 'Rudagatan 18'
Set-TextValue 84 6 '64.72822387894983'
Set-TextValue 84 7 '21.066092511399347'

# Row 85
$ws.Cells.Item(85, 1).Value = 'Individ 15'
$ws.Cells.Item(85, 2).Value = 5
$ws.Cells.Item(85, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(85, 4).Value = 'smoker,dog,>18'
$ws.Cells.Item(85, 5).Value = 'Getargatan 13'
Set-TextValue 85 6 '64.72515988986189'
Set-TextValue 85 7 '21.081590156187488'

# Row 86
$ws.Cells.Item(86, 1).Value = 'Individ 18'
$ws.Cells.Item(86, 2).Value = 20
$ws.Cells.Item(86, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(86, 4).Value = 'license,>18'
$ws.Cells.Item(86, 5).Value = 'OmvÃ¤gen 11'
Set-TextValue 86 6 '64.722700477568'
Set-TextValue 86 7 '21.076038097904636'

# Row 87
$ws.Cells.Item(87, 1).Value = 'Individ 20'
$ws.Cells.Item(87, 2).Value = 5
$ws.Cells.Item(87, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(87, 4).Value = 'cat'
$ws.Cells.Item(87, 5).Value = 'BergsbyvÃ¤gen 19'
Set-TextValue 87 6 '64.72253627756791'
Set-TextValue 87 7 '21.081495540234208'

# Row 88
$ws.Cells.Item(88, 1).Value = 'Individ 21'
$ws.Cells.Item(88, 2).Value = 5
$ws.Cells.Item(88, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(88, 4).Value = '>18'
$ws.Cells.Item(88, 5).Value = 'Nybyggargatan 12'
Set-TextValue 88 6 '64.72279765876095'
Set-TextValue 88 7 '21.096477997904646'

# Row 89
$ws.Cells.Item(89, 1).Value = 'Individ 24'
$ws.Cells.Item(89, 2).Value = 20
$ws.Cells.Item(89, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(89, 4).Value = ''
$ws.Cells.Item(89, 5).Value = 'Roddargatan 12'
Set-TextValue 89 6 '64.71749269547423'
Set-TextValue 89 7 '21.09246261139871'

# Row 90
$ws.Cells.Item(90, 1).Value = 'Individ 28'
$ws.Cells.Item(90, 2).Value = 15
$ws.Cells.Item(90, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(90, 4).Value = 'license,cat'
$ws.Cells.Item(90, 5).Value = 'Lillgatan 7A'
Set-TextValue 90 6 '64.71565457756343'
Set-TextValue 90 7 '21.0969892402338'

# Row 91
$ws.Cells.Item(91, 1).Value = 'Individ 29'
$ws.Cells.Item(91, 2).Value = 15
$ws.Cells.Item(91, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(91, 4).Value = 'license'
$ws.Cells.Item(91, 5).Value = 'Roddargatan 29'
Set-TextValue 91 6 '64.7158097128137'
Set-TextValue 91 7 '21.098717311398616'

# Row 92
$ws.Cells.Item(92, 1).Value = 'Individ 41'
$ws.Cells.Item(92, 2).Value = 10
$ws.Cells.Item(92, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(92, 4).Value = 'license'
$ws.Cells.Item(92, 5).Value = 'SkelleftehamnsvÃ¤gen 113'
Set-TextValue 92 6 '64.7145646253258'
Set-TextValue 92 7 '21.16073439790411'

# Row 93
$ws.Cells.Item(93, 1).Value = 'Individ 44'
$ws.Cells.Item(93, 2).Value = 10
$ws.Cells.Item(93, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(93, 4).Value = ''
$ws.Cells.Item(93, 5).Value = 'VÃ¤nskapsgatan 4'
Set-TextValue 93 6 '64.71183961148486'
Set-TextValue 93 7 '21.17002340955187'

# Row 94
$ws.Cells.Item(94, 1).Value = 'Individ 59'
$ws.Cells.Item(94, 2).Value = 10
$ws.Cells.Item(94, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(94, 4).Value = ''
$ws.Cells.Item(94, 5).Value = 'Lotsens grÃ¤nd 5'
Set-TextValue 94 6 '64.69594034903153'
Set-TextValue 94 7 '21.190897069067788'

# Row 95
$ws.Cells.Item(95, 1).Value = 'Individ 62'
$ws.Cells.Item(95, 2).Value = 60
$ws.Cells.Item(95, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(95, 4).Value = ''
$ws.Cells.Item(95, 5).Value = 'VÃ¥gens grÃ¤nd 2'
Set-TextValue 95 6 '64.69386766364026'
Set-TextValue 95 7 '21.195669282562054'

# Row 96
$ws.Cells.Item(96, 1).Value = 'Individ 69'
$ws.Cells.Item(96, 2).Value = 45
$ws.Cells.Item(96, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(96, 4).Value = ''
$ws.Cells.Item(96, 5).Value = 'BokgrÃ¤nd 7'
Set-TextValue 96 6 '64.71491237756301'
Set-TextValue 96 7 '21.158029382563345'

# Row 97
$ws.Cells.Item(97, 1).Value = 'Individ 1'
$ws.Cells.Item(97, 2).Value = 30
$ws.Cells.Item(97, 3).Value = '(''Förmiddag'', ''9-11'')'
$ws.Cells.Item(97, 4).Value = 'license,dog,woman'
$ws.Cells.Item(97, 5).Value = 'Rudagatan 51'
Set-TextValue 97 6 '64.73076905004987'
Set-TextValue 97 7 '21.062869127380175'

# Row 98
$ws.Cells.Item(98, 1).Value = 'Individ 28'
$ws.Cells.Item(98, 2).Value = 20
$ws.Cells.Item(98, 3).Value = '(''Förmiddag'', ''9-11'')'
$ws.Cells.Item(98, 4).Value = 'license,cat'
$ws.Cells.Item(98, 5).Value = 'Lillgatan 7A'
Set-TextValue 98 6 '64.71565457756343'
Set-TextValue 98 7 '21.0969892402338'

# Row 99
$ws.Cells.Item(99, 1).Value = 'Individ 29'
$ws.Cells.Item(99, 2).Value = 40
$ws.Cells.Item(99, 3).Value = '(''Förmiddag'', ''9-11'')'
$ws.Cells.Item(99, 4).Value = 'license'
$ws.Cells.Item(99, 5).Value = 'Roddargatan 29'
Set-TextValue 99 6 '64.7158097128137'
Set-TextValue 99 7 '21.098717311398616'

# Row 100
$ws.Cells.Item(100, 1).Value = 'Individ 37'
$ws.Cells.Item(100, 2).Value = 20
$ws.Cells.Item(100, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(100, 4).Value = 'license'
$ws.Cells.Item(100, 5).Value = 'MÃ¥bÃ¤rsgatan 12'
Set-TextValue 100 6 '64.71651135982128'
Set-TextValue 100 7 '21.15559144023381'

# Row 101
$ws.Cells.Item(101, 1).Value = 'Individ 38'
$ws.Cells.Item(101, 2).Value = 45
$ws.Cells.Item(101, 3).Value = '(''Eftermiddag'', ''13-15'')'
$ws.Cells.Item(101, 4).Value = 'license'
$ws.Cells.Item(101, 5).Value = 'Karagangatan 30'
Set-TextValue 101 6 '64.71678347084324'
Set-TextValue 101 7 '21.15778882489306'

